$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.821.54"
$ws.Range("E2").Value = "  -0.87%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.810.67"
$ws.Range("E3").Value = "  +0.58%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"
$ws.Range("E4").Value = "  +0.16%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "308.82"
$ws.Range("E5").Value = "  +0.24%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.002"
$ws.Range("E6").Value = "  +0.07%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4315"
$ws.Range("E7").Value = "  +2.42%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3699"
$ws.Range("E8").Value = "  +2.83%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07251"
$ws.Range("E9").Value = "  -0.65%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8681"

$ws.Range("B11").Value = "WrappedEther"
$ws.Range("C11").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.991.20"
$ws.Range("E11").Value = "  +6.51%  "

$ws.Range("B12").Value = "Solana"
$ws.Range("C12").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "20.85"
$ws.Range("E12").Value = "  +2.85%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.629"
$ws.Range("E13").Value = "  +3.84%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.352"
$ws.Range("E14").Value = "  +0.98%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.06906"
$ws.Range("E15").Value = "  +1.95%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.003"
$ws.Range("E16").Value = "  +0.04%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "80.46"
$ws.Range("E17").Value = "  -0.08%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008848"
$ws.Range("E18").Value = "  +0.89%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.001"
$ws.Range("E19").Value = "  -0.05%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.27"
$ws.Range("E20").Value = "  +1.75%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "26.865.50"
$ws.Range("E21").Value = "  -1.61%  "

$ws.Range("E22").Value = "  +2.43%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.18"
$ws.Range("E23").Value = "  +1.23%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.206.43"
$ws.Range("E24").Value = "  +5.85%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "153.58"
$ws.Range("E25").Value = "  +0.09%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.871"
$ws.Range("E26").Value = "  -2.99%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.26"
$ws.Range("E27").Value = "  +0.55%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.211"
$ws.Range("E28").Value = "  +3.47%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.906"
$ws.Range("E29").Value = "  +14.94%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "115.35"
$ws.Range("E30").Value = "  +1.65%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08943"
$ws.Range("E31").Value = "  -0.83%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7571"
$ws.Range("E32").Value = "  +3.35%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.169"
$ws.Range("E33").Value = "  +6.48%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.441"
$ws.Range("E34").Value = "  +2.09%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.780"
$ws.Range("E35").Value = "  -2.83%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.005"
$ws.Range("E36").Value = "  +0.36%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.125"
$ws.Range("E37").Value = "  +3.88%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05218"
$ws.Range("E38").Value = "  +1.21%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01925"
$ws.Range("E39").Value = "  +0.88%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5082"
$ws.Range("E40").Value = "  +1.84%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1648"
$ws.Range("E41").Value = "  +0.88%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.664"
$ws.Range("E42").Value = "  +0.72%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.544"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.288"
$ws.Range("E44").Value = "  +2.53%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "106.48"
$ws.Range("E45").Value = "  +0.97%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.40"
$ws.Range("E46").Value = "  +0.86%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.002"
$ws.Range("E47").Value = "  +0.07%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.657"
$ws.Range("E48").Value = "  +3.25%  "

$ws.Range("B49").Value = "Decentraland"
$ws.Range("C49").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.4571"
$ws.Range("E49").Value = "  +0.39%  "

$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06279"
$ws.Range("E50").Value = "  -0.44%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.801"
$ws.Range("E51").Value = "  +4.00%  "
